$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row of data (row 39) continuing the daily user impact status log
$ws.Range("A39").Value = 45992
$ws.Range("B39").Value = 5612
$ws.Range("C39").Value = 4205
$ws.Range("D39").Value = 3873
$ws.Range("E39").Value = 245
$ws.Range("F39").Value = 51
$ws.Range("G39").Value = 36
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0

# Match date number format used by the row above (A38) for the new date cell
$ws.Range("A39").NumberFormat = $ws.Range("A38").NumberFormat

# Update the selected cell / sqref to reflect the newly added last row
$ws.Range("A39:I39").Select()
